# Auto-generated edit script.
# Applies literal value corrections to the H:N (market-price / profit) columns
# across several worksheets, as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 46503.637  # H17: 62313.715 -> 46503.637
$ws.Cells.Item(17, 9).Value = 0  # I17: 500 -> 0
$ws.Cells.Item(17, 10).Value = 46503.637  # J17: 64131.766 -> 46503.637
$ws.Cells.Item(17, 11).Value = 0  # K17: 1500 -> 0
$ws.Cells.Item(17, 12).Value = 139510.911  # L17: 192395.298 -> 139510.911
$ws.Cells.Item(17, 13).ClearContents()  # M17: remove cell (was None)
$ws.Cells.Item(17, 14).Value = -139846.911  # N17: -192731.298 -> -139846.911
$ws.Cells.Item(88, 8).Value = 5950.654  # H88: 6848.048 -> 5950.654
$ws.Cells.Item(88, 9).Value = 6050.909  # I88: 5321.3 -> 6050.909
$ws.Cells.Item(88, 10).Value = 5877.1333  # J88: 8236 -> 5877.1333
$ws.Cells.Item(88, 11).Value = 6050.909  # K88: 5321.3 -> 6050.909
$ws.Cells.Item(88, 12).Value = 5877.1333  # L88: 8236 -> 5877.1333
$ws.Cells.Item(88, 13).Value = -5644.909  # M88: -4915.3 -> -5644.909
$ws.Cells.Item(88, 14).Value = -6689.1333  # N88: -9048 -> -6689.1333
$ws.Cells.Item(91, 8).Value = 5950.654  # H91: 6848.048 -> 5950.654
$ws.Cells.Item(91, 9).Value = 6050.909  # I91: 5321.3 -> 6050.909
$ws.Cells.Item(91, 10).Value = 5877.1333  # J91: 8236 -> 5877.1333
$ws.Cells.Item(91, 11).Value = 6050.909  # K91: 5321.3 -> 6050.909
$ws.Cells.Item(91, 12).Value = 5877.1333  # L91: 8236 -> 5877.1333
$ws.Cells.Item(91, 13).Value = -4646.909  # M91: -3917.3 -> -4646.909
$ws.Cells.Item(91, 14).Value = -8685.133300000001  # N91: -11044 -> -8685.133300000001
$ws.Cells.Item(98, 8).Value = 1418.3143  # H98: 1608.4445 -> 1418.3143
$ws.Cells.Item(98, 9).Value = 645.96295  # I98: 1256.7 -> 645.96295
$ws.Cells.Item(98, 10).Value = 4025  # J98: 2613.4285 -> 4025
$ws.Cells.Item(98, 11).Value = 645.96295  # K98: 1256.7 -> 645.96295
$ws.Cells.Item(98, 12).Value = 4025  # L98: 2613.4285 -> 4025
$ws.Cells.Item(98, 13).Value = 852.03705  # M98: 241.3 -> 852.03705
$ws.Cells.Item(98, 14).Value = -7021  # N98: -5609.4285 -> -7021
$ws.Cells.Item(122, 8).Value = 1418.3143  # H122: 1608.4445 -> 1418.3143
$ws.Cells.Item(122, 9).Value = 645.96295  # I122: 1256.7 -> 645.96295
$ws.Cells.Item(122, 10).Value = 4025  # J122: 2613.4285 -> 4025
$ws.Cells.Item(122, 11).Value = 1937.88885  # K122: 3770.1 -> 1937.88885
$ws.Cells.Item(122, 12).Value = 12075  # L122: 7840.2855 -> 12075
$ws.Cells.Item(122, 13).Value = 512.1111500000002  # M122: -1320.1 -> 512.1111500000002
$ws.Cells.Item(122, 14).Value = -16975  # N122: -12740.2855 -> -16975
$ws.Cells.Item(132, 8).Value = 3392233.2  # H132: 2859043.5 -> 3392233.2
$ws.Cells.Item(132, 9).Value = 4083858.8  # I132: 3176352 -> 4083858.8
$ws.Cells.Item(132, 10).Value = 3268.2  # J132: 3267.4285 -> 3268.2
$ws.Cells.Item(132, 11).Value = 12251576.4  # K132: 9529056 -> 12251576.4
$ws.Cells.Item(132, 12).Value = 9804.599999999999  # L132: 9802.2855 -> 9804.599999999999
$ws.Cells.Item(132, 13).Value = -12249046.4  # M132: -9526526 -> -12249046.4
$ws.Cells.Item(132, 14).Value = -14864.6  # N132: -14862.2855 -> -14864.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3504.36  # H32: 1717.21 -> 3504.36
$ws.Cells.Item(32, 9).Value = 2839.956  # I32: 1705.4375 -> 2839.956
$ws.Cells.Item(32, 10).Value = 10222.223  # J32: 1999.75 -> 10222.223
$ws.Cells.Item(32, 11).Value = 2839.956  # K32: 1705.4375 -> 2839.956
$ws.Cells.Item(32, 12).Value = 10222.223  # L32: 1999.75 -> 10222.223
$ws.Cells.Item(32, 13).Value = -2552.956  # M32: -1418.4375 -> -2552.956
$ws.Cells.Item(32, 14).Value = -10796.223  # N32: -2573.75 -> -10796.223
$ws.Cells.Item(61, 8).Value = 1420.8206  # H61: 1241.8334 -> 1420.8206
$ws.Cells.Item(61, 9).Value = 747.41174  # I61: 665.3022999999999 -> 747.41174
$ws.Cells.Item(61, 10).Value = 6000  # J61: 6200 -> 6000
$ws.Cells.Item(61, 11).Value = 747.41174  # K61: 665.3022999999999 -> 747.41174
$ws.Cells.Item(61, 12).Value = 6000  # L61: 6200 -> 6000
$ws.Cells.Item(61, 13).Value = -535.41174  # M61: -453.3022999999999 -> -535.41174
$ws.Cells.Item(61, 14).Value = -6424  # N61: -6624 -> -6424
$ws.Cells.Item(63, 8).Value = 3221  # H63: 3337.9546 -> 3221
$ws.Cells.Item(63, 9).Value = 2049.9285  # I63: 2068.2144 -> 2049.9285
$ws.Cells.Item(63, 10).Value = 6500  # J63: 5560 -> 6500
$ws.Cells.Item(63, 11).Value = 2049.9285  # K63: 2068.2144 -> 2049.9285
$ws.Cells.Item(63, 12).Value = 6500  # L63: 5560 -> 6500
$ws.Cells.Item(63, 13).Value = -1363.9285  # M63: -1382.2144 -> -1363.9285
$ws.Cells.Item(63, 14).Value = -7872  # N63: -6932 -> -7872
$ws.Cells.Item(66, 8).Value = 3221  # H66: 3337.9546 -> 3221
$ws.Cells.Item(66, 9).Value = 2049.9285  # I66: 2068.2144 -> 2049.9285
$ws.Cells.Item(66, 10).Value = 6500  # J66: 5560 -> 6500
$ws.Cells.Item(66, 11).Value = 10249.6425  # K66: 10341.072 -> 10249.6425
$ws.Cells.Item(66, 12).Value = 32500  # L66: 27800 -> 32500
$ws.Cells.Item(66, 13).Value = -6817.6425  # M66: -6909.072 -> -6817.6425
$ws.Cells.Item(66, 14).Value = -39364  # N66: -34664 -> -39364
$ws.Cells.Item(110, 8).Value = 2107.7856  # H110: 968.63635 -> 2107.7856
$ws.Cells.Item(110, 9).Value = 501  # I110: 344.85715 -> 501
$ws.Cells.Item(110, 10).Value = 5000  # J110: 2060.25 -> 5000
$ws.Cells.Item(110, 11).Value = 501  # K110: 344.85715 -> 501
$ws.Cells.Item(110, 12).Value = 5000  # L110: 2060.25 -> 5000
$ws.Cells.Item(110, 13).Value = 1544  # M110: 1700.14285 -> 1544
$ws.Cells.Item(110, 14).Value = -9090  # N110: -6150.25 -> -9090
$ws.Cells.Item(132, 8).Value = 2785.7297  # H132: 2923.7646 -> 2785.7297
$ws.Cells.Item(132, 9).Value = 2230.3572  # I132: 2304.8462 -> 2230.3572
$ws.Cells.Item(132, 10).Value = 4513.5557  # J132: 4935.25 -> 4513.5557
$ws.Cells.Item(132, 11).Value = 6691.071599999999  # K132: 6914.5386 -> 6691.071599999999
$ws.Cells.Item(132, 12).Value = 13540.6671  # L132: 14805.75 -> 13540.6671
$ws.Cells.Item(132, 13).Value = -4161.071599999999  # M132: -4384.5386 -> -4161.071599999999
$ws.Cells.Item(132, 14).Value = -18600.6671  # N132: -19865.75 -> -18600.6671
$ws.Cells.Item(136, 8).Value = 1420.8206  # H136: 1241.8334 -> 1420.8206
$ws.Cells.Item(136, 9).Value = 747.41174  # I136: 665.3022999999999 -> 747.41174
$ws.Cells.Item(136, 10).Value = 6000  # J136: 6200 -> 6000
$ws.Cells.Item(136, 11).Value = 2242.23522  # K136: 1995.9069 -> 2242.23522
$ws.Cells.Item(136, 12).Value = 18000  # L136: 18600 -> 18000
$ws.Cells.Item(136, 13).Value = 307.76478  # M136: 554.0931 -> 307.76478
$ws.Cells.Item(136, 14).Value = -23100  # N136: -23700 -> -23100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2502.547  # H107: 3124.8538 -> 2502.547
$ws.Cells.Item(107, 9).Value = 2099.658  # I107: 2661.9355 -> 2099.658
$ws.Cells.Item(107, 10).Value = 3523.2  # J107: 4559.9 -> 3523.2
$ws.Cells.Item(107, 11).Value = 2099.658  # K107: 2661.9355 -> 2099.658
$ws.Cells.Item(107, 12).Value = 3523.2  # L107: 4559.9 -> 3523.2
$ws.Cells.Item(107, 13).Value = -179.6579999999999  # M107: -741.9355 -> -179.6579999999999
$ws.Cells.Item(107, 14).Value = -7363.2  # N107: -8399.9 -> -7363.2
$ws.Cells.Item(134, 8).Value = 2606.1555  # H134: 2681.7727 -> 2606.1555
$ws.Cells.Item(134, 9).Value = 2052.7297  # I134: 2135.9143 -> 2052.7297
$ws.Cells.Item(134, 10).Value = 5165.75  # J134: 4804.5557 -> 5165.75
$ws.Cells.Item(134, 11).Value = 6158.1891  # K134: 6407.742899999999 -> 6158.1891
$ws.Cells.Item(134, 12).Value = 15497.25  # L134: 14413.6671 -> 15497.25
$ws.Cells.Item(134, 13).Value = -3623.1891  # M134: -3872.742899999999 -> -3623.1891
$ws.Cells.Item(134, 14).Value = -20567.25  # N134: -19483.6671 -> -20567.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2465.3809  # H31: 2523.4033 -> 2465.3809
$ws.Cells.Item(31, 9).Value = 1470.6383  # I31: 1533.3636 -> 1470.6383
$ws.Cells.Item(31, 10).Value = 5387.4375  # J31: 4943.5 -> 5387.4375
$ws.Cells.Item(31, 11).Value = 1470.6383  # K31: 1533.3636 -> 1470.6383
$ws.Cells.Item(31, 12).Value = 5387.4375  # L31: 4943.5 -> 5387.4375
$ws.Cells.Item(31, 13).Value = -1175.6383  # M31: -1238.3636 -> -1175.6383
$ws.Cells.Item(31, 14).Value = -5977.4375  # N31: -5533.5 -> -5977.4375
$ws.Cells.Item(34, 8).Value = 2465.3809  # H34: 2523.4033 -> 2465.3809
$ws.Cells.Item(34, 9).Value = 1470.6383  # I34: 1533.3636 -> 1470.6383
$ws.Cells.Item(34, 10).Value = 5387.4375  # J34: 4943.5 -> 5387.4375
$ws.Cells.Item(34, 11).Value = 1470.6383  # K34: 1533.3636 -> 1470.6383
$ws.Cells.Item(34, 12).Value = 5387.4375  # L34: 4943.5 -> 5387.4375
$ws.Cells.Item(34, 13).Value = -1268.6383  # M34: -1331.3636 -> -1268.6383
$ws.Cells.Item(34, 14).Value = -5791.4375  # N34: -5347.5 -> -5791.4375
$ws.Cells.Item(86, 8).Value = 6583.5835  # H86: 7731.231 -> 6583.5835
$ws.Cells.Item(86, 9).Value = 4667.3335  # I86: 6084.5 -> 4667.3335
$ws.Cells.Item(86, 10).Value = 8499.833000000001  # J86: 9142.714 -> 8499.833000000001
$ws.Cells.Item(86, 11).Value = 4667.3335  # K86: 6084.5 -> 4667.3335
$ws.Cells.Item(86, 12).Value = 8499.833000000001  # L86: 9142.714 -> 8499.833000000001
$ws.Cells.Item(86, 13).Value = -3544.3335  # M86: -4961.5 -> -3544.3335
$ws.Cells.Item(86, 14).Value = -10745.833  # N86: -11388.714 -> -10745.833
$ws.Cells.Item(89, 8).Value = 6583.5835  # H89: 7731.231 -> 6583.5835
$ws.Cells.Item(89, 9).Value = 4667.3335  # I89: 6084.5 -> 4667.3335
$ws.Cells.Item(89, 10).Value = 8499.833000000001  # J89: 9142.714 -> 8499.833000000001
$ws.Cells.Item(89, 11).Value = 23336.6675  # K89: 30422.5 -> 23336.6675
$ws.Cells.Item(89, 12).Value = 42499.165  # L89: 45713.57 -> 42499.165
$ws.Cells.Item(89, 13).Value = -17720.6675  # M89: -24806.5 -> -17720.6675
$ws.Cells.Item(89, 14).Value = -53731.165  # N89: -56945.57 -> -53731.165
$ws.Cells.Item(99, 8).Value = 2472.5  # H99: 2184.6667 -> 2472.5
$ws.Cells.Item(99, 9).Value = 1575  # I99: 1544.4445 -> 1575
$ws.Cells.Item(99, 10).Value = 4267.5  # J99: 3145 -> 4267.5
$ws.Cells.Item(99, 11).Value = 1575  # K99: 1544.4445 -> 1575
$ws.Cells.Item(99, 12).Value = 4267.5  # L99: 3145 -> 4267.5
$ws.Cells.Item(99, 13).Value = -77  # M99: -46.44450000000006 -> -77
$ws.Cells.Item(99, 14).Value = -7263.5  # N99: -6141 -> -7263.5
$ws.Cells.Item(126, 8).Value = 2472.5  # H126: 2184.6667 -> 2472.5
$ws.Cells.Item(126, 9).Value = 1575  # I126: 1544.4445 -> 1575
$ws.Cells.Item(126, 10).Value = 4267.5  # J126: 3145 -> 4267.5
$ws.Cells.Item(126, 11).Value = 4725  # K126: 4633.333500000001 -> 4725
$ws.Cells.Item(126, 12).Value = 12802.5  # L126: 9435 -> 12802.5
$ws.Cells.Item(126, 13).Value = -2255  # M126: -2163.333500000001 -> -2255
$ws.Cells.Item(126, 14).Value = -17742.5  # N126: -14375 -> -17742.5
$ws.Cells.Item(132, 8).Value = 1841.6666  # H132: 1874.5 -> 1841.6666
$ws.Cells.Item(132, 9).Value = 1341.9286  # I132: 1398.0952 -> 1341.9286
$ws.Cells.Item(132, 10).Value = 4173.778  # J132: 3875.4 -> 4173.778
$ws.Cells.Item(132, 11).Value = 4025.7858  # K132: 4194.2856 -> 4025.7858
$ws.Cells.Item(132, 12).Value = 12521.334  # L132: 11626.2 -> 12521.334
$ws.Cells.Item(132, 13).Value = -1495.7858  # M132: -1664.2856 -> -1495.7858
$ws.Cells.Item(132, 14).Value = -17581.334  # N132: -16686.2 -> -17581.334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 998.6667  # H5: 913.1429000000001 -> 998.6667
$ws.Cells.Item(5, 10).Value = 2300  # J5: 1666.6666 -> 2300
$ws.Cells.Item(5, 12).Value = 6900  # L5: 4999.9998 -> 6900
$ws.Cells.Item(5, 14).Value = -7124  # N5: -5223.9998 -> -7124
$ws.Cells.Item(122, 8).Value = 1131.4286  # H122: 1537.5385 -> 1131.4286
$ws.Cells.Item(122, 9).Value = 712  # I122: 930 -> 712
$ws.Cells.Item(122, 10).Value = 1341.1428  # J122: 1719.8 -> 1341.1428
$ws.Cells.Item(122, 11).Value = 6408  # K122: 8370 -> 6408
$ws.Cells.Item(122, 12).Value = 12070.2852  # L122: 15478.2 -> 12070.2852
$ws.Cells.Item(122, 13).Value = -3958  # M122: -5920 -> -3958
$ws.Cells.Item(122, 14).Value = -16970.2852  # N122: -20378.2 -> -16970.2852
$ws.Cells.Item(131, 8).Value = 2085.4167  # H131: 1826.9445 -> 2085.4167
$ws.Cells.Item(131, 9).Value = 4021.6667  # I131: 4016.6667 -> 4021.6667
$ws.Cells.Item(131, 10).Value = 1440  # J131: 1389 -> 1440
$ws.Cells.Item(131, 11).Value = 12065.0001  # K131: 12050.0001 -> 12065.0001
$ws.Cells.Item(131, 12).Value = 4320  # L131: 4167 -> 4320
$ws.Cells.Item(131, 13).Value = -7025.000100000001  # M131: -7010.000100000001 -> -7025.000100000001
$ws.Cells.Item(131, 14).Value = -14400  # N131: -14247 -> -14400
$ws.Cells.Item(132, 8).Value = 1509.5  # H132: 1201.7097 -> 1509.5
$ws.Cells.Item(132, 9).Value = 906  # I132: 816.2381 -> 906
$ws.Cells.Item(132, 10).Value = 3320  # J132: 2011.2 -> 3320
$ws.Cells.Item(132, 11).Value = 8154  # K132: 7346.142900000001 -> 8154
$ws.Cells.Item(132, 12).Value = 29880  # L132: 18100.8 -> 29880
$ws.Cells.Item(132, 13).Value = -5624  # M132: -4816.142900000001 -> -5624
$ws.Cells.Item(132, 14).Value = -34940  # N132: -23160.8 -> -34940
$ws.Cells.Item(133, 8).Value = 3923.5652  # H133: 4227.905 -> 3923.5652
$ws.Cells.Item(133, 9).Value = 4267.9  # I133: 4813.3335 -> 4267.9
$ws.Cells.Item(133, 10).Value = 3658.6924  # J133: 3788.8333 -> 3658.6924
$ws.Cells.Item(133, 11).Value = 12803.7  # K133: 14440.0005 -> 12803.7
$ws.Cells.Item(133, 12).Value = 10976.0772  # L133: 11366.4999 -> 10976.0772
$ws.Cells.Item(133, 13).Value = -7743.699999999999  # M133: -9380.000499999998 -> -7743.699999999999
$ws.Cells.Item(133, 14).Value = -21096.0772  # N133: -21486.4999 -> -21096.0772
$ws.Cells.Item(134, 8).Value = 2941.6365  # H134: 2495.2727 -> 2941.6365
$ws.Cells.Item(134, 9).Value = 1908.2858  # I134: 1249.2307 -> 1908.2858
$ws.Cells.Item(134, 10).Value = 4750  # J134: 4295.1113 -> 4750
$ws.Cells.Item(134, 11).Value = 5724.857400000001  # K134: 3747.6921 -> 5724.857400000001
$ws.Cells.Item(134, 12).Value = 14250  # L134: 12885.3339 -> 14250
$ws.Cells.Item(134, 13).Value = -654.8574000000008  # M134: 1322.3079 -> -654.8574000000008
$ws.Cells.Item(134, 14).Value = -24390  # N134: -23025.3339 -> -24390
$ws.Cells.Item(135, 8).Value = 998.6667  # H135: 913.1429000000001 -> 998.6667
$ws.Cells.Item(135, 10).Value = 2300  # J135: 1666.6666 -> 2300
$ws.Cells.Item(135, 12).Value = 20700  # L135: 14999.9994 -> 20700
$ws.Cells.Item(135, 14).Value = -25770  # N135: -20069.9994 -> -25770
$ws.Cells.Item(136, 8).Value = 1687.6  # H136: 1779.1428 -> 1687.6
$ws.Cells.Item(136, 9).Value = 1214.6666  # I136: 1315.2941 -> 1214.6666
$ws.Cells.Item(136, 10).Value = 2397  # J136: 2496 -> 2397
$ws.Cells.Item(136, 11).Value = 3643.9998  # K136: 3945.8823 -> 3643.9998
$ws.Cells.Item(136, 12).Value = 7191  # L136: 7488 -> 7191
$ws.Cells.Item(136, 13).Value = 1456.0002  # M136: 1154.1177 -> 1456.0002
$ws.Cells.Item(136, 14).Value = -17391  # N136: -17688 -> -17391
$ws.Cells.Item(137, 8).Value = 2768.8064  # H137: 2968.16 -> 2768.8064
$ws.Cells.Item(137, 9).Value = 2270.5881  # I137: 2321.6667 -> 2270.5881
$ws.Cells.Item(137, 10).Value = 3373.7856  # J137: 3564.923 -> 3373.7856
$ws.Cells.Item(137, 11).Value = 6811.7643  # K137: 6965.000100000001 -> 6811.7643
$ws.Cells.Item(137, 12).Value = 10121.3568  # L137: 10694.769 -> 10121.3568
$ws.Cells.Item(137, 13).Value = -1711.7643  # M137: -1865.000100000001 -> -1711.7643
$ws.Cells.Item(137, 14).Value = -20321.3568  # N137: -20894.769 -> -20321.3568
$ws.Cells.Item(138, 8).Value = 3068.5334  # H138: 2809.3635 -> 3068.5334
$ws.Cells.Item(138, 9).Value = 1247.2858  # I138: 874.2 -> 1247.2858
$ws.Cells.Item(138, 10).Value = 4662.125  # J138: 4422 -> 4662.125
$ws.Cells.Item(138, 11).Value = 3741.8574  # K138: 2622.6 -> 3741.8574
$ws.Cells.Item(138, 12).Value = 13986.375  # L138: 13266 -> 13986.375
$ws.Cells.Item(138, 13).Value = 1398.1426  # M138: 2517.4 -> 1398.1426
$ws.Cells.Item(138, 14).Value = -24266.375  # N138: -23546 -> -24266.375
$ws.Cells.Item(139, 8).Value = 12503684  # H139: 3478.7856 -> 12503684
$ws.Cells.Item(139, 9).Value = 27781426  # I139: 2777.8572 -> 27781426
$ws.Cells.Item(139, 10).Value = 3713.2727  # J139: 4179.7144 -> 3713.2727
$ws.Cells.Item(139, 11).Value = 83344278  # K139: 8333.571599999999 -> 83344278
$ws.Cells.Item(139, 12).Value = 11139.8181  # L139: 12539.1432 -> 11139.8181
$ws.Cells.Item(139, 13).Value = -83339138  # M139: -3193.571599999999 -> -83339138
$ws.Cells.Item(139, 14).Value = -21419.8181  # N139: -22819.1432 -> -21419.8181
$ws.Cells.Item(140, 8).Value = 7578694.5  # H140: 4216.6665 -> 7578694.5
$ws.Cells.Item(140, 9).Value = 15152501  # I140: 1500 -> 15152501
$ws.Cells.Item(140, 10).Value = 4888.1816  # J140: 4760 -> 4888.1816
$ws.Cells.Item(140, 11).Value = 45457503  # K140: 4500 -> 45457503
$ws.Cells.Item(140, 12).Value = 14664.5448  # L140: 14280 -> 14664.5448
$ws.Cells.Item(140, 13).Value = -45452323  # M140: 680 -> -45452323
$ws.Cells.Item(140, 14).Value = -25024.5448  # N140: -24640 -> -25024.5448
$ws.Cells.Item(141, 8).Value = 2857.1428  # H141: 3260 -> 2857.1428
$ws.Cells.Item(141, 9).Value = 2750  # I141: 3433.3333 -> 2750
$ws.Cells.Item(141, 11).Value = 8250  # K141: 10299.9999 -> 8250
$ws.Cells.Item(141, 13).Value = -3070  # M141: -5119.999899999999 -> -3070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1345.8422  # H113: 1789.2084 -> 1345.8422
$ws.Cells.Item(113, 9).Value = 865.05554  # I113: 1067.2858 -> 865.05554
$ws.Cells.Item(113, 10).Value = 10000  # J113: 2799.9 -> 10000
$ws.Cells.Item(113, 11).Value = 865.05554  # K113: 1067.2858 -> 865.05554
$ws.Cells.Item(113, 12).Value = 10000  # L113: 2799.9 -> 10000
$ws.Cells.Item(113, 13).Value = 1304.94446  # M113: 1102.7142 -> 1304.94446
$ws.Cells.Item(113, 14).Value = -14340  # N113: -7139.9 -> -14340

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1429.7  # H7: 1402.9354 -> 1429.7
$ws.Cells.Item(7, 10).Value = 2589.889  # J7: 2390.9 -> 2589.889
$ws.Cells.Item(7, 12).Value = 2589.889  # L7: 2390.9 -> 2589.889
$ws.Cells.Item(7, 14).Value = -2813.889  # N7: -2614.9 -> -2813.889
$ws.Cells.Item(40, 8).Value = 3286.6086  # H40: 6519.5 -> 3286.6086
$ws.Cells.Item(40, 9).Value = 3288.4443  # I40: 14400 -> 3288.4443
$ws.Cells.Item(40, 10).Value = 3280  # J40: 3142.1428 -> 3280
$ws.Cells.Item(40, 11).Value = 3288.4443  # K40: 14400 -> 3288.4443
$ws.Cells.Item(40, 12).Value = 3280  # L40: 3142.1428 -> 3280
$ws.Cells.Item(40, 13).Value = -3152.4443  # M40: -14264 -> -3152.4443
$ws.Cells.Item(40, 14).Value = -3552  # N40: -3414.1428 -> -3552
$ws.Cells.Item(126, 8).Value = 1429.7  # H126: 1402.9354 -> 1429.7
$ws.Cells.Item(126, 10).Value = 2589.889  # J126: 2390.9 -> 2589.889
$ws.Cells.Item(126, 12).Value = 7769.667  # L126: 7172.700000000001 -> 7769.667
$ws.Cells.Item(126, 14).Value = -12709.667  # N126: -12112.7 -> -12709.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 296149.5  # H122: 502000.94 -> 296149.5
$ws.Cells.Item(122, 9).Value = 372066.22  # I122: 589830.5600000001 -> 372066.22
$ws.Cells.Item(122, 10).Value = 3327.8572  # J122: 4299.6665 -> 3327.8572
$ws.Cells.Item(122, 11).Value = 1116198.66  # K122: 1769491.68 -> 1116198.66
$ws.Cells.Item(122, 12).Value = 9983.571599999999  # L122: 12898.9995 -> 9983.571599999999
$ws.Cells.Item(122, 13).Value = -1113748.66  # M122: -1767041.68 -> -1113748.66
$ws.Cells.Item(122, 14).Value = -14883.5716  # N122: -17798.9995 -> -14883.5716
$ws.Cells.Item(132, 8).Value = 14275.733  # H132: 8162.4053 -> 14275.733
$ws.Cells.Item(132, 9).Value = 3895.7727  # I132: 1545.2413 -> 3895.7727
$ws.Cells.Item(132, 10).Value = 24204.39  # J132: 26438.38 -> 24204.39
$ws.Cells.Item(132, 11).Value = 11687.3181  # K132: 4635.7239 -> 11687.3181
$ws.Cells.Item(132, 12).Value = 72613.17  # L132: 79315.14 -> 72613.17
$ws.Cells.Item(132, 13).Value = -9157.3181  # M132: -2105.7239 -> -9157.3181
$ws.Cells.Item(132, 14).Value = -77673.17  # N132: -84375.14 -> -77673.17
